# Updated cryptos list (Price / Volume(1h) refresh, plus a Polygon/WrappedEther
# row-order swap in rows 12-13) to mirror the latest GitHub Actions scrape.
#
# Note: Price values are stored as text in this sheet (e.g. "29.940.74",
# "1.002"). Where the new value would otherwise be auto-recognised by Excel
# as a number (breaking the trailing zeros / text formatting of the source
# data), we prefix it with a leading apostrophe so Excel keeps it as text,
# exactly like the original cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.940.74'
$ws.Range('E2').Value = '  +0.45%  '

$ws.Range('D3').Value = '1.891.86'
$ws.Range('E3').Value = '  -0.46%  '

$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  +0.22%  '

$ws.Range('D5').Value = '''0.8188'
$ws.Range('E5').Value = '  +6.48%  '

$ws.Range('D6').Value = '''241.79'
$ws.Range('E6').Value = '  +0.57%  '

$ws.Range('D7').Value = '''1.002'
$ws.Range('E7').Value = '  +0.26%  '

$ws.Range('D8').Value = '''0.3228'
$ws.Range('E8').Value = '  +5.63%  '

$ws.Range('D9').Value = '''26.42'
$ws.Range('E9').Value = '  +3.63%  '

$ws.Range('D10').Value = '''0.07029'
$ws.Range('E10').Value = '  +2.60%  '

$ws.Range('D11').Value = '''0.08037'
$ws.Range('E11').Value = '  +0.66%  '

$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = '''0.7453'
$ws.Range('E12').Value = '  +0.94%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.899.71'
$ws.Range('E13').Value = '  -0.32%  '

$ws.Range('E14').Value = '  +0.27%  '

$ws.Range('D15').Value = '''92.28'
$ws.Range('E15').Value = '  +1.14%  '

$ws.Range('D16').Value = '29.940.90'
$ws.Range('E16').Value = '  +0.37%  '

$ws.Range('D17').Value = '''14.02'
$ws.Range('E17').Value = '  +1.95%  '

$ws.Range('D18').Value = '''5.890'

$ws.Range('D19').Value = '''244.55'
$ws.Range('E19').Value = '  -0.20%  '

$ws.Range('D20').Value = '''0.000007750'
$ws.Range('E20').Value = '  +0.43%  '

$ws.Range('E21').Value = '  +0.31%  '

$ws.Range('D22').Value = '2.153.69'
$ws.Range('E22').Value = '  -0.83%  '

$ws.Range('D23').Value = '''1.002'
$ws.Range('E23').Value = '  +0.24%  '

$ws.Range('D24').Value = '''6.905'
$ws.Range('E24').Value = '  -0.36%  '

$ws.Range('D25').Value = '''0.1575'
$ws.Range('E25').Value = '  +22.35%  '

$ws.Range('D26').Value = '''166.01'
$ws.Range('E26').Value = '  -0.51%  '

$ws.Range('D27').Value = '''9.184'
$ws.Range('E27').Value = '  -0.74%  '

$ws.Range('D28').Value = '''18.86'
$ws.Range('E28').Value = '  +0.89%  '

$ws.Range('D29').Value = '''2.071'
$ws.Range('E29').Value = '  +1.76%  '

$ws.Range('E30').Value = '  -2.19%  '

$ws.Range('D31').Value = '''1.522'
$ws.Range('E31').Value = '  +0.68%  '

$ws.Range('D32').Value = '''4.268'
$ws.Range('E32').Value = '  -0.22%  '

$ws.Range('D33').Value = '''0.05589'
$ws.Range('E33').Value = '  +6.13%  '

$ws.Range('E34').Value = '  +0.10%  '

$ws.Range('D35').Value = '''1.270'
$ws.Range('E35').Value = '  +2.08%  '

$ws.Range('D36').Value = '''0.7311'
$ws.Range('E36').Value = '  +0.59%  '

$ws.Range('D37').Value = '''2.722'
$ws.Range('E37').Value = '  +0.17%  '

$ws.Range('D38').Value = '''0.01914'
$ws.Range('E38').Value = '  +0.26%  '

$ws.Range('D39').Value = '''2.786'
$ws.Range('E39').Value = '  +0.21%  '

$ws.Range('D40').Value = '''0.4413'
$ws.Range('E40').Value = '  -0.04%  '

$ws.Range('D41').Value = '''71.93'
$ws.Range('E41').Value = '  -0.40%  '

$ws.Range('D42').Value = '''5.950'
$ws.Range('E42').Value = '  -4.14%  '

$ws.Range('E43').Value = '  +1.02%  '

$ws.Range('D44').Value = '''1.002'
$ws.Range('E44').Value = '  +0.25%  '

$ws.Range('E45').Value = '  -0.31%  '

$ws.Range('D46').Value = '''100.73'
$ws.Range('E46').Value = '  +0.70%  '

$ws.Range('D47').Value = '''7.568'
$ws.Range('E47').Value = '  -0.50%  '

$ws.Range('D48').Value = '''9.700'
$ws.Range('E48').Value = '  -0.55%  '

$ws.Range('D49').Value = '''990.85'
$ws.Range('E49').Value = '  +8.68%  '

$ws.Range('D50').Value = '2.051.37'
$ws.Range('E50').Value = '  -0.47%  '

$ws.Range('D51').Value = '''36.02'
$ws.Range('E51').Value = '  -0.53%  '
